{"js": "// The \"Requisitos\" bullet list paragraph is reshuffled into a new order\n// (one run per \"LOxxxxx - Name (Requisito)\" entry, each followed by a\n// manual line break). This rebuilds that paragraph with the runs in the\n// new order while leaving every other paragraph in the document alone.\n\nconst NEW_ORDER = [\n  \"LOB1053 -  F\u00edsica III  (Requisito)\",\n  \"LOB1045 -  Leitura e Produ\u00e7\u00e3o de Textos Acad\u00eamicos  (Requisito)\",\n  \"LOB1008 -  Ci\u00eancia, Tecnologia e Sociedade  (Requisito)\",\n  \"LOB1036 -  Geometria Anal\u00edtica  (Requisito)\",\n  \"LOB1037 -  \u00c0lgebra Linear  (Requisito)\",\n  \"LOB1041 -  F\u00edsica Experimental II  (Requisito)\",\n  \"LOB1042 -  F\u00edsica Experimental IV  (Requisito)\",\n  \"LOQ4095 -  Qu\u00edmica Geral Experimental  (Requisito)\",\n  \"LOB1039 -  F\u00edsica Experimental III  (Requisito)\",\n  \"LOB1018 -  F\u00edsica I  (Requisito)\",\n  \"LOQ4100 -  Fundamentos de Qu\u00edmica para Engenharia I (Requisito)\",\n  \"LOB1004 -  C\u00e1lculo II  (Requisito)\",\n  \"LOB1038 -  F\u00edsica Experimental I  (Requisito)\",\n  \"LOB1052 -  C\u00e1lculo III  (Requisito)\",\n  \"LOM3236 -  Processos de Fabrica\u00e7\u00e3o  (Requisito)\",\n  \"LOM3261 -  M\u00e9todos Num\u00e9ricos e Aplica\u00e7\u00f5es  (Requisito)\",\n  \"LOM3218 -  Introdu\u00e7\u00e3o \u00e0 Engenharia F\u00edsica  (Requisito)\",\n  \"LOB1003 -  C\u00e1lculo I  (Requisito)\",\n  \"LOB1006 -  C\u00e1lculo IV  (Requisito)\",\n  \"LOM3241 -  Qu\u00edmica de Materiais  (Requisito)\",\n  \"LOB1021 -  F\u00edsica IV  (Requisito)\",\n  \"LOM3016 -  Introdu\u00e7\u00e3o \u00e0  Ci\u00eancia dos Materiais  (Requisito)\",\n  \"LOM3260 -  Computa\u00e7\u00e3o Cient\u00edfica em Python  (Requisito)\",\n  \"LOM3204 -  Desenho T\u00e9cnico e Projeto Assistido por Computador  (Requisito)\",\n  \"LOB1012 -  Estat\u00edstica  (Requisito)\",\n  \"LOB1019 -  F\u00edsica II  (Requisito)\",\n];\n\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text,items/style\");\nawait context.sync();\n\n// Locate the \"Requisitos\" Heading2 paragraph, then the bullet-list\n// paragraph that immediately follows it (the one this edit targets).\nlet reqIndex = -1;\nfor (let i = 0; i < paras.items.length; i++) {\n  if (paras.items[i].text.trim() === \"Requisitos\") {\n    reqIndex = i;\n    break;\n  }\n}\nif (reqIndex === -1 || reqIndex + 1 >= paras.items.length) {\n  throw new Error(\"Could not locate the 'Requisitos' list paragraph\");\n}\nconst listPara = paras.items[reqIndex + 1];\n\n// Unique marker text (very unlikely to collide with real content) used to\n// split \"freshly inserted new content\" from \"old content still pending\n// removal\" without relying on fragile character offsets.\nconst SENTINEL = \"\\u0001REORDER-SENTINEL\\u0001\";\n\n// Mark the original start of the paragraph's content.\nconst startMarker = listPara.getRange(\"Start\");\nstartMarker.insertText(SENTINEL, \"Before\");\nawait context.sync();\n\n// Insert every entry of the new order, in reverse, always \"Before\" the\n// (fixed) paragraph start -- each call lands just ahead of the previous\n// one, so the visible order ends up matching NEW_ORDER top to bottom.\nconst insertionPoint = listPara.getRange(\"Start\");\nfor (let i = NEW_ORDER.length - 1; i >= 0; i--) {\n  insertionPoint.insertText(NEW_ORDER[i] + \"\\u000b\", \"Before\");\n}\nawait context.sync();\n\n// Find the sentinel we planted earlier, then delete everything from the\n// sentinel to the end of the paragraph -- that span is exactly the\n// original (now-superseded) run content plus the sentinel itself.\nconst found = listPara.search(SENTINEL, { matchCase: true });\nfound.load(\"items\");\nawait context.sync();\nif (found.items.length === 0) {\n  throw new Error(\"Lost track of the reorder sentinel\");\n}\nconst tail = found.items[0].expandTo(listPara.getRange(\"End\"));\ntail.delete();\nawait context.sync();\n", "ps1": "# The \"Requisitos\" bullet list paragraph is reshuffled into a new order\n# (one run per \"LOxxxxx - Name (Requisito)\" entry, each followed by a\n# manual line break). This rebuilds that paragraph with the runs in the\n# new order while leaving every other paragraph in the document alone.\n\n$NewOrder = @(\n    \"LOB1053 -  F\u00edsica III  (Requisito)\",\n    \"LOB1045 -  Leitura e Produ\u00e7\u00e3o de Textos Acad\u00eamicos  (Requisito)\",\n    \"LOB1008 -  Ci\u00eancia, Tecnologia e Sociedade  (Requisito)\",\n    \"LOB1036 -  Geometria Anal\u00edtica  (Requisito)\",\n    \"LOB1037 -  \u00c0lgebra Linear  (Requisito)\",\n    \"LOB1041 -  F\u00edsica Experimental II  (Requisito)\",\n    \"LOB1042 -  F\u00edsica Experimental IV  (Requisito)\",\n    \"LOQ4095 -  Qu\u00edmica Geral Experimental  (Requisito)\",\n    \"LOB1039 -  F\u00edsica Experimental III  (Requisito)\",\n    \"LOB1018 -  F\u00edsica I  (Requisito)\",\n    \"LOQ4100 -  Fundamentos de Qu\u00edmica para Engenharia I (Requisito)\",\n    \"LOB1004 -  C\u00e1lculo II  (Requisito)\",\n    \"LOB1038 -  F\u00edsica Experimental I  (Requisito)\",\n    \"LOB1052 -  C\u00e1lculo III  (Requisito)\",\n    \"LOM3236 -  Processos de Fabrica\u00e7\u00e3o  (Requisito)\",\n    \"LOM3261 -  M\u00e9todos Num\u00e9ricos e Aplica\u00e7\u00f5es  (Requisito)\",\n    \"LOM3218 -  Introdu\u00e7\u00e3o \u00e0 Engenharia F\u00edsica  (Requisito)\",\n    \"LOB1003 -  C\u00e1lculo I  (Requisito)\",\n    \"LOB1006 -  C\u00e1lculo IV  (Requisito)\",\n    \"LOM3241 -  Qu\u00edmica de Materiais  (Requisito)\",\n    \"LOB1021 -  F\u00edsica IV  (Requisito)\",\n    \"LOM3016 -  Introdu\u00e7\u00e3o \u00e0  Ci\u00eancia dos Materiais  (Requisito)\",\n    \"LOM3260 -  Computa\u00e7\u00e3o Cient\u00edfica em Python  (Requisito)\",\n    \"LOM3204 -  Desenho T\u00e9cnico e Projeto Assistido por Computador  (Requisito)\",\n    \"LOB1012 -  Estat\u00edstica  (Requisito)\",\n    \"LOB1019 -  F\u00edsica II  (Requisito)\"\n)\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n$count = $paras.Count\n\n# Locate the \"Requisitos\" Heading2 paragraph, then the bullet-list\n# paragraph that immediately follows it (the one this edit targets).\n$listParaIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    if ($paras.Item($i).Range.Text.Trim() -eq \"Requisitos\") {\n        $listParaIndex = $i + 1\n        break\n    }\n}\nif ($listParaIndex -eq -1) {\n    throw \"Could not locate the 'Requisitos' list paragraph\"\n}\n$listPara = $paras.Item($listParaIndex)\n$listRange = $listPara.Range\n\n# Unique marker text (very unlikely to collide with real content) used to\n# split \"freshly inserted new content\" from \"old content still pending\n# removal\" without relying on fragile character offsets.\n$Sentinel = [char]1 + \"REORDER-SENTINEL\" + [char]1\n\n$paraStart = $listRange.Start\n$marker = $d.Range($paraStart, $paraStart)\n$marker.InsertBefore($Sentinel)\n\n# Insert every entry of the new order, in reverse, always \"Before\" the\n# (fixed) paragraph start -- each call lands just ahead of the previous\n# one, so the visible order ends up matching $NewOrder top to bottom.\nfor ($i = $NewOrder.Count - 1; $i -ge 0; $i--) {\n    $ins = $d.Range($paraStart, $paraStart)\n    $ins.InsertBefore($NewOrder[$i] + [char]11)\n}\n\n# Find the sentinel we planted earlier, then delete everything from the\n# sentinel to the end of the paragraph -- that span is exactly the\n# original (now-superseded) run content plus the sentinel itself.\n$searchRange = $listPara.Range\n$searchRange.Find.ClearFormatting()\n$searchRange.Find.Forward = $true\n$searchRange.Find.Wrap = 0\n$found = $searchRange.Find.Execute($Sentinel)\nif (-not $found) {\n    throw \"Lost track of the reorder sentinel\"\n}\n\n$tailStart = $searchRange.Start\n$paraEnd = $listPara.Range.End\n$tail = $d.Range($tailStart, $paraEnd)\n$tail.Delete()\n"}
